$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): drop the bold/centered/bordered style and the
#     "Unnamed: 0" label in A1 -------------------------------------------
$ws.Range("A1:AM1").ClearFormats()
$ws.Range("A1").ClearContents()

# --- Row 3 (Revisit count) ------------------------------------------------
$ws.Range("B3").Value = 5
$ws.Range("E3").Value = 14
$ws.Range("H3").Value = 13
$ws.Range("I3").Value = 24
$ws.Range("J3").Value = 31
$ws.Range("K3").Value = 29
$ws.Range("M3").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 27
$ws.Range("Z3").Value = 1
$ws.Range("AC3").Value = 25
$ws.Range("AM3").ClearContents()

# --- Row 4 (Fixation count) ------------------------------------------------
$ws.Range("B4").Value = 7
$ws.Range("E4").Value = 20
$ws.Range("H4").Value = 27
$ws.Range("I4").Value = 146
$ws.Range("J4").Value = 65
$ws.Range("K4").Value = 64
$ws.Range("M4").Value = 2
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 134
$ws.Range("Z4").Value = 2
$ws.Range("AC4").Value = 151
$ws.Range("AM4").ClearContents()

# --- Row 5 (Dwell time (ms)) ------------------------------------------------
$ws.Range("B5").Value = 4771.73
$ws.Range("E5").Value = 10562.7
$ws.Range("H5").Value = 13672.63
$ws.Range("I5").Value = 57712.02
$ws.Range("J5").Value = 31910.26
$ws.Range("K5").Value = 28263.76
$ws.Range("M5").Value = 316.97
$ws.Range("R5").Value = 1484.82
$ws.Range("S5").Value = 55526.38
$ws.Range("Z5").Value = 316.97
$ws.Range("AC5").Value = 70264.37
$ws.Range("AM5").ClearContents()

# --- Row 6 (Dwell time (%)) ------------------------------------------------
$ws.Range("B6").Value = 2.68
$ws.Range("C6").Value = 3.11
$ws.Range("E6").Value = 5.93
$ws.Range("F6").Value = 0.79
$ws.Range("G6").Value = 3.22
$ws.Range("H6").Value = 7.68
$ws.Range("I6").Value = 32.41
$ws.Range("J6").Value = 17.92
$ws.Range("K6").Value = 15.87
$ws.Range("L6").Value = 4.67
$ws.Range("M6").Value = 0.18
$ws.Range("R6").Value = 0.83
$ws.Range("S6").Value = 31.19
$ws.Range("T6").Value = 0.37
$ws.Range("U6").Value = 0.12
$ws.Range("W6").Value = 0.52
$ws.Range("X6").Value = 10.41
$ws.Range("Y6").Value = 7.15
$ws.Range("Z6").Value = 0.18
$ws.Range("AA6").Value = 0.41
$ws.Range("AC6").Value = 39.46
$ws.Range("AD6").Value = 0.09
$ws.Range("AE6").Value = 0.12
$ws.Range("AF6").Value = 1.32
$ws.Range("AG6").Value = 4.21
$ws.Range("AH6").Value = 3.95
$ws.Range("AL6").Value = 0.23
$ws.Range("AM6").ClearContents()

# --- Row 7 (Fixation duration (ms)) ------------------------------------------------
$ws.Range("B7").Value = 681.6799999999999
$ws.Range("E7").Value = 528.14
$ws.Range("H7").Value = 506.39
$ws.Range("I7").Value = 395.29
$ws.Range("J7").Value = 490.93
$ws.Range("K7").Value = 441.62
$ws.Range("M7").Value = 158.48
$ws.Range("R7").Value = 742.41
$ws.Range("S7").Value = 414.38
$ws.Range("Z7").Value = 158.48
$ws.Range("AC7").Value = 465.33
$ws.Range("AM7").ClearContents()

# --- Row 8 (First fixation duration (ms)) -----------------------------------
$ws.Range("AM8").ClearContents()

# --- Drop the trailing blank row 11 -----------------------------------------
$ws.Rows.Item(11).Delete()

Write-Host "edit complete"
